# Fixed annotation and simplified repeated code
# Appends 4 new trial rows (12-15) to the HFP_TEST sheet, following the
# same column layout as the existing rows (A:U).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("TEST", 1, 1, 2023, 10, 30, 15, 23, 30, 200, 25, 125, 65, 65, 65, 113, 21, 43, 27, 51, 128),
    @("TEST", 1, 2, 2023, 10, 30, 15, 23, 40, 250, 250, 250, 250, 250, 250, 250, 250, 250, 250, 250, 128),
    @("TEST", 1, 3, 2023, 10, 30, 15, 23, 54, 35, 75, 0, 0, 0, 14, 14, 14, 14, 14, 14, 128),
    @("TEST", 8, 1, 2023, 10, 30, 15, 43, 16, 52, 52, 50, 30, 0, 10, 10, 10, 16, 16, 13, 255)
)

$startRow = 12
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowValues = $newRows[$i]
    $rowIndex = $startRow + $i
    for ($col = 1; $col -le $rowValues.Count; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowValues[$col - 1]
    }
}
